$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header format (from G1, the "sum" header) into H1 so the
# new "Save" header matches the other header cells' style (bold, centered,
# bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add numeric values for the new Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
